$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on column D cells we touch so numeric-looking strings
# are preserved exactly as text (matching original inlineStr cell type).
$dCells = @('D2', 'D3', 'D4', 'D5', 'D6', 'D7', 'D8', 'D9', 'D10', 'D11', 'D12', 'D13', 'D14', 'D15', 'D16', 'D17', 'D18', 'D19', 'D20', 'D22', 'D23', 'D28', 'D40', 'D41', 'D42', 'D43', 'D44', 'D45', 'D47', 'D48', 'D49', 'D50')
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values
$ws.Range('D2').Value = '246.33'
$ws.Range('D3').Value = '24.09'
$ws.Range('D4').Value = '5.359'
$ws.Range('D5').Value = '0.05797'
$ws.Range('D6').Value = '6.467'
$ws.Range('D7').Value = '3.328'
$ws.Range('D8').Value = '0.8091'
$ws.Range('D9').Value = '0.9221'
$ws.Range('B10').Value = 'One'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D10').Value = '0.01067'
$ws.Range('E10').Value = '9OneONEBestin24h'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D11').Value = '0.1411'
$ws.Range('E11').Value = '10WazirXWRX'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').Value = '0.07358'
$ws.Range('E12').Value = '11MandalaExchangeTokenMDX'
$ws.Range('B13').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C13').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D13').Value = '0.03166'
$ws.Range('E13').Value = '12LiechtensteinCryptoassetsExchangeLCX'
$ws.Range('B14').Value = 'BitrueCoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D14').Value = '0.03059'
$ws.Range('E14').Value = '13BitrueCoinBTR'
$ws.Range('B15').Value = 'BitMartToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D15').Value = '0.09358'
$ws.Range('E15').Value = '14BitMartTokenBMX'
$ws.Range('B16').Value = 'MCDex'
$ws.Range('C16').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D16').Value = '3.855'
$ws.Range('E16').Value = '15MCDexMCB'
$ws.Range('B17').Value = 'BitForexToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D17').Value = '0.001559'
$ws.Range('E17').Value = '16BitForexTokenBF'
$ws.Range('B18').Value = 'CoinExToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D18').Value = '0.04738'
$ws.Range('E18').Value = '17CoinExTokenCET'
$ws.Range('D19').Value = '0.005896'
$ws.Range('D20').Value = '0.001277'
$ws.Range('D22').Value = '0.00008805'
$ws.Range('E22').Value = '21NitroExNTX'
$ws.Range('D23').Value = '3.608'
$ws.Range('D28').Value = '0.0002351'
$ws.Range('D40').Value = '0.03826'
$ws.Range('D41').Value = '0.006434'
$ws.Range('B42').Value = 'CEJI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range('D42').Value = '0.004102'
$ws.Range('E42').Value = '41CEJICEJI'
$ws.Range('B43').Value = 'BKEXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('D43').Value = '0.1064'
$ws.Range('E43').Value = '42BKEXTokenBKK'
$ws.Range('D44').Value = '0.008424'
$ws.Range('D45').Value = '0.00005331'
$ws.Range('D47').Value = '0.6859'
$ws.Range('D48').Value = '0.001845'
$ws.Range('E48').Value = '47BOLOBOLOWorstin24h'
$ws.Range('D49').Value = '0.00002101'
$ws.Range('D50').Value = '0.0002001'
